$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename G1 header
$ws.Range("G1").Value = "new_regression"

# Move/replace G2:G11 with the text values that used to be in H (all TRUE),
# then delete the H column entirely.
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "'TRUE"
    $cell.Style = "Normal"
}

# Delete column H so the table collapses back to A:G
$ws.Range("H1").EntireColumn.Delete()
